# Edit for commit: "added 230 packet run 2d"
# Target sheet: "CNN 2D" (sheet2.xml)
# 1. Fill in the previously-empty Acc/Loss/Time data for the "230 Packets"
#    column group (Q3:S51) -- the header labels (Q2/R2/S2, R1) already existed.
# 2. Lay down the header placeholders for the *next* column group, "210 Packets"
#    (V1 title + U2/V2/W2 sub-headers), without its data yet (matches the diff).
# 3. Refresh the view: scroll so column C is left-most, select U11.
# 4. Set the page to portrait orientation (adds <pageSetup .../> on save).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CNN 2D")
$ws.Activate()

# --- 1. "230 Packets" run data (Q3:S51) ---------------------------------
$rowNums = @(3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,48,49,50,51)
$qVals   = @(88.412737846374498,82.941639423370304,84.075224399566594,82.941639423370304,83.226221799850407,86.332917213439899,88.690209388732896,87.734484672546301,87.900489568710299,85.071265697479205,87.957406044006305,87.917089462280202,87.615907192230196,87.855434417724595,87.703657150268498,87.276780605316105,82.941639423370304,83.283138275146399,87.060970067977905,88.296532630920396,87.914717197418199,86.958998441696096,88.059383630752507,88.104444742202702,83.271276950836096,87.760573625564504,82.941639423370304,87.264925241470294,87.226980924606295,87.430930137634206,87.032514810562105,87.136858701705904,87.848317623138399,88.607203960418701,87.162947654724107,87.051486968994098,87.630134820938096,88.021439313888493,82.939267158508301,87.727367877960205,88.547915220260606,87.563735246658297,87.784284353256197,84.599334001541095,87.696540355682302,86.5155220031738,87.554246187209998,83.828586339950505,82.941639423370304)
$rVals   = @(0.264765702028815,2.6234215325207999,0.29043662896042799,2.6234215325207999,0.30168453050423799,0.29237627452062398,0.27188473447171702,0.263791751192403,0.26024197134928401,0.31405787178449801,0.265459077894315,0.26089043518615401,0.28684454690901301,0.26828533962800399,0.27761200636785699,0.26761227144483102,2.6234215325207999,0.30064781244692101,0.29719339798510602,0.25886780820845401,0.27213141955012099,0.29416847245012201,0.27134637639290798,0.26532468423970001,0.30881494839443802,0.27090088435611598,0.29699016838172398,0.284858503280708,0.29508772968735603,0.30049246508136601,0.28769832397700501,0.291208854473755,0.26365141086413701,0.25171889708067702,0.27219134221932501,0.28005267312842902,0.29075864017798397,0.25811993370596498,0.33297670761921799,0.26665605488314398,0.272359962511594,0.27786745282299602,0.28105299032230302,0.29336210926296502,0.26635324594339999,0.28276195560040801,0.27109047767019301,0.32420824389931802,0.336088248761089)
$sVals   = @(80.053917646407996,72.853944778442298,73.184900760650606,72.305456638336096,73.291011810302706,73.248565196990896,73.021291494369507,73.007450342178302,69.657107353210407,72.425482749938894,75.334001779556203,71.083654880523596,72.312381505966101,75.232361078262301,71.242110967636094,70.763823270797701,72.717964172363196,73.800097465515094,72.217993497848497,75.353200435638399,71.697496414184499,75.1012024879455,72.344587802886906,71.357617855071993,73.223740577697697,69.036902189254704,73.070634126663194,72.2075066566467,73.820480346679602,72.832002878189002,73.749711275100694,70.985456943511906,75.029891729354802,75.250156164169297,73.761211633682194,70.9909183979034,73.186587095260606,74.147364854812594,72.172680377960205,72.596275568008394,74.676091432571397,72.845376014709402,72.744814634323106,72.313711166381793,71.462985992431598,71.436714410781804,71.229358673095703,71.258295536041203,69.337912321090698)

for ($i = 0; $i -lt $rowNums.Length; $i++) {
  $r = $rowNums[$i]
  $ws.Cells.Item($r, 17).Value = $qVals[$i]
  $ws.Cells.Item($r, 18).Value = $rVals[$i]
  $ws.Cells.Item($r, 19).Value = $sVals[$i]
}

# --- 2. Header placeholders for the next column group, "210 Packets" ---
$ws.Range("V1").Value = "210 Packets"
$ws.Range("U2").Value = "Acc"
$ws.Range("V2").Value = "Loss"
$ws.Range("W2").Value = "Time"

# --- 3. View state: scroll to column C, select U11 ----------------------
$ws.Range("U11").Select()
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1

# --- 4. Page orientation --------------------------------------------------
$ws.PageSetup.Orientation = 1
